$wb = $excel.ActiveWorkbook

# --- Sheet 1: Forecast Comparison (numeric MyForecast values) ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 178
$ws1.Range("D3").Value = 179
$ws1.Range("D4").Value = 76
$ws1.Range("D14").Value = 116
$ws1.Range("D16").Value = 107

# --- Sheet 2: Summary (these values are stored as text in the workbook,
#     so a leading apostrophe is used to force Excel to keep the
#     numeric-looking / date-looking strings as literal text instead of
#     auto-converting them to a number or date serial) ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'1919"
$ws2.Range("B10").Value = "'1054"
$ws2.Range("B11").Value = "'557"
$ws2.Range("B12").Value = "'179"
$ws2.Range("B14").Value = "'60"
$ws2.Range("B15").Value = "'2025-04-27"
